# Refresh the cryptos list "Price" (D) and "Volume(1h)" (E) columns,
# per the scraper's GitHub Actions commit. Each value is written with a
# leading apostrophe so Excel stores it as literal text (these look like
# numbers/dates, e.g. "239.37" or "29.792.39", but the source workbook
# keeps them as plain strings). Resetting Style to "Normal" afterwards
# clears the quote-prefix formatting that the apostrophe trick applies,
# so the cell keeps the same (default) style it started with.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Ref = "D2"; Value = '29.792.39'},
    @{Ref = "E2"; Value = '  -1.44%  '},
    @{Ref = "D3"; Value = '1.892.44'},
    @{Ref = "E3"; Value = '  -1.10%  '},
    @{Ref = "E4"; Value = '  -0.04%  '},
    @{Ref = "D5"; Value = '0.7569'},
    @{Ref = "E5"; Value = '  +2.48%  '},
    @{Ref = "D6"; Value = '239.37'},
    @{Ref = "E6"; Value = '  -1.51%  '},
    @{Ref = "D8"; Value = '0.3041'},
    @{Ref = "E8"; Value = '  -2.82%  '},
    @{Ref = "D9"; Value = '25.19'},
    @{Ref = "E9"; Value = '  -6.49%  '},
    @{Ref = "D10"; Value = '0.06824'},
    @{Ref = "E10"; Value = '  -1.76%  '},
    @{Ref = "D11"; Value = '0.07964'},
    @{Ref = "E11"; Value = '  -0.18%  '},
    @{Ref = "D12"; Value = '0.7476'},
    @{Ref = "E12"; Value = '  -3.31%  '},
    @{Ref = "D13"; Value = '1.898.37'},
    @{Ref = "E13"; Value = '  -1.42%  '},
    @{Ref = "D14"; Value = '5.204'},
    @{Ref = "E14"; Value = '  -1.01%  '},
    @{Ref = "D15"; Value = '90.99'},
    @{Ref = "E15"; Value = '  -0.44%  '},
    @{Ref = "D16"; Value = '29.798.21'},
    @{Ref = "E16"; Value = '  -1.60%  '},
    @{Ref = "D17"; Value = '6.017'},
    @{Ref = "E17"; Value = '  +3.82%  '},
    @{Ref = "D18"; Value = '13.84'},
    @{Ref = "E18"; Value = '  -2.64%  '},
    @{Ref = "E19"; Value = '  -1.82%  '},
    @{Ref = "D20"; Value = '233.47'},
    @{Ref = "E20"; Value = '  -4.72%  '},
    @{Ref = "E21"; Value = '  -0.04%  '},
    @{Ref = "D22"; Value = '2.146.21'},
    @{Ref = "E22"; Value = '  -2.74%  '},
    @{Ref = "E23"; Value = '  -0.05%  '},
    @{Ref = "D24"; Value = '6.944'},
    @{Ref = "E24"; Value = '  +4.92%  '},
    @{Ref = "D25"; Value = '9.235'},
    @{Ref = "E25"; Value = '  -1.59%  '},
    @{Ref = "D26"; Value = '164.96'},
    @{Ref = "E26"; Value = '  -0.24%  '},
    @{Ref = "E27"; Value = '  -1.62%  '},
    @{Ref = "D28"; Value = '0.1289'},
    @{Ref = "E28"; Value = '  +1.46%  '},
    @{Ref = "D29"; Value = '2.048'},
    @{Ref = "E29"; Value = '  -3.77%  '},
    @{Ref = "D30"; Value = '1.339'},
    @{Ref = "E30"; Value = '  -0.91%  '},
    @{Ref = "D31"; Value = '1.514'},
    @{Ref = "E31"; Value = '  -1.92%  '},
    @{Ref = "D32"; Value = '4.273'},
    @{Ref = "E32"; Value = '  -0.91%  '},
    @{Ref = "D33"; Value = '4.003'},
    @{Ref = "E33"; Value = '  -1.83%  '},
    @{Ref = "D34"; Value = '0.05277'},
    @{Ref = "E34"; Value = '  +2.16%  '},
    @{Ref = "D35"; Value = '1.239'},
    @{Ref = "E35"; Value = '  -4.20%  '},
    @{Ref = "D36"; Value = '0.7288'},
    @{Ref = "E36"; Value = '  -2.76%  '},
    @{Ref = "D37"; Value = '2.711'},
    @{Ref = "E37"; Value = '  -1.88%  '},
    @{Ref = "D38"; Value = '0.01922'},
    @{Ref = "E38"; Value = '  -0.77%  '},
    @{Ref = "D39"; Value = '2.760'},
    @{Ref = "E39"; Value = '  -0.67%  '},
    @{Ref = "D40"; Value = '6.201'},
    @{Ref = "E40"; Value = '  -2.89%  '},
    @{Ref = "D41"; Value = '0.4403'},
    @{Ref = "E41"; Value = '  -1.48%  '},
    @{Ref = "D42"; Value = '72.42'},
    @{Ref = "E42"; Value = '  -4.58%  '},
    @{Ref = "E43"; Value = '  -1.71%  '},
    @{Ref = "E44"; Value = '  +0.01%  '},
    @{Ref = "D45"; Value = '0.8246'},
    @{Ref = "E45"; Value = '  -0.98%  '},
    @{Ref = "D46"; Value = '101.05'},
    @{Ref = "E46"; Value = '  -0.27%  '},
    @{Ref = "D47"; Value = '7.584'},
    @{Ref = "E47"; Value = '  -0.83%  '},
    @{Ref = "D48"; Value = '9.799'},
    @{Ref = "E48"; Value = '  -0.52%  '},
    @{Ref = "D49"; Value = '2.051.64'},
    @{Ref = "E49"; Value = '  -1.91%  '},
    @{Ref = "D50"; Value = '35.94'},
    @{Ref = "E50"; Value = '  -2.60%  '},
    @{Ref = "D51"; Value = '0.05943'},
    @{Ref = "E51"; Value = '  -0.45%  '}
)

foreach ($u in $updates) {
    $ws.Range($u.Ref).Value = "'" + $u.Value
    $ws.Range($u.Ref).Style = "Normal"
}
